$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table3 "Overview": add row 9 for the new file
# d452f1dd-74d8-4c8f-972c-7be67665c439.md
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = "d452f1dd-74d8-4c8f-972c-7be67665c439.md"
$wsOverview.Range("B9").Value = "e2e\d452f1dd-74d8-4c8f-972c-7be67665c439.md"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-21 20:55:32"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d452f1dd74d84c8f972c7be67665c4396a0279bf/e2e/d452f1dd-74d8-4c8f-972c-7be67665c439.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\d452f1dd-74d8-4c8f-972c-7be67665c439.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table1 "zh-cn": add row 9
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A9").Value = "d452f1dd-74d8-4c8f-972c-7be67665c439.md"
$wsZh.Range("B9").Value = ".md"
$wsZh.Range("C9").Value = "Ready for handoff"
$wsZh.Range("D9").Value = "e2e"
$wsZh.Range("E9").Value = "ht"
$wsZh.Range("F9").Value = "False"
$wsZh.Range("G9").Value = "d452f1dd-74d8-4c8f-972c-7be67665c439.b07018d703e746de6f5096f3692d55169f1d988e.zh-cn.xlf"
$wsZh.Range("H9").Value = "2016-08-21 20:55:28"
$wsZh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I9").Value = ""
$wsZh.Range("J9").Value = ""
$wsZh.Range("K9").Value = "0001-01-01 00:00:00"
$wsZh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L9").Value = ""
$wsZh.Range("M9").Value = "True"
$wsZh.Range("N9").Value = ""
$wsZh.Range("O9").Value = "False"
$wsZh.Range("P9").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d452f1dd74d84c8f972c7be67665c4396a0279bf/e2e/d452f1dd-74d8-4c8f-972c-7be67665c439.md", [System.Type]::Missing, [System.Type]::Missing, "d452f1dd-74d8-4c8f-972c-7be67665c439.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table2 "de-de": add row 9
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A9").Value = "d452f1dd-74d8-4c8f-972c-7be67665c439.md"
$wsDe.Range("B9").Value = ".md"
$wsDe.Range("C9").Value = "Ready for handoff"
$wsDe.Range("D9").Value = "e2e"
$wsDe.Range("E9").Value = "ht"
$wsDe.Range("F9").Value = "False"
$wsDe.Range("G9").Value = "d452f1dd-74d8-4c8f-972c-7be67665c439.b07018d703e746de6f5096f3692d55169f1d988e.de-de.xlf"
$wsDe.Range("H9").Value = "2016-08-21 20:55:32"
$wsDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I9").Value = ""
$wsDe.Range("J9").Value = ""
$wsDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L9").Value = ""
$wsDe.Range("M9").Value = "True"
$wsDe.Range("N9").Value = ""
$wsDe.Range("O9").Value = "False"
$wsDe.Range("P9").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d452f1dd74d84c8f972c7be67665c4396a0279bf/e2e/d452f1dd-74d8-4c8f-972c-7be67665c439.md", [System.Type]::Missing, [System.Type]::Missing, "d452f1dd-74d8-4c8f-972c-7be67665c439.md") | Out-Null
